# Swap the presentation's applied theme colours from the "Integral" theme
# back to the default "Office" theme colours (Design tab -> Themes -> Office).
#
# The underlying edit recorded in the source OOXML swaps the contents of
# ppt/theme/theme1.xml ("Office Theme") and ppt/theme/theme2.xml ("Integral")
# so that the slide master (which is wired to theme2.xml) ends up showing the
# standard Office colour palette. The font scheme (major/minor Latin = Arial)
# and format scheme are identical between the two themes, so only the colour
# scheme actually needs to change.
#
# We reproduce that visible effect with the supported COM surface by editing
# the active theme's ThemeColorScheme in place, one slot at a time, using the
# 12 standard Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) in the canonical clrScheme order.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# RGB() isn't available in this host, so colours are passed as the packed
# 0x00BBGGRR integer that OLE_COLOR / the .RGB property uses.
# index : clrScheme slot : target hex (RRGGBB) : packed BGR integer
$themeColors.Item(1).RGB  = 0          # dk1      000000
$themeColors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Item(3).RGB  = 6968388    # dk2      44546A
$themeColors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Item(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Item(8).RGB  = 49407      # accent4  FFC000
$themeColors.Item(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Item(10).RGB = 4697456    # accent6  70AD47
$themeColors.Item(11).RGB = 12673797   # hlink    0563C1
$themeColors.Item(12).RGB = 7491477    # folHlink 954F72
